# Add a new "localdb" command-category to the hidden '#system' lookup sheet
# and register it as a named range, mirroring how every other category
# (aws.s3, macro, mail, ...) is laid out: a header cell holding the
# category name, followed by the list of its function signatures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1. Make room for the new category by inserting a new column at N.
#    Everything that used to live in N..AC (macro, mail, number, pdf,
#    rdbms, redis, sms, sound, ssh, step, web, webalert, webcookie, ws,
#    ws.async, xml) shifts right to O..AD.
$ws.Range("N1").EntireColumn.Insert()

# 2. Populate the new column with the "localdb" category header and its
#    six function signatures.
$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# 3. Column A ("target") holds the sorted list of every category name
#    (used to populate a picker), independent of the unrelated lists
#    that happen to share the same row numbers in other columns
#    (E, G, X, ...). Shift only the A14:A29 values down to A15:A30 so
#    "localdb" can be inserted in its alphabetical spot - between
#    "json" and "macro" - without disturbing any other column.
for ($r = 29; $r -ge 14; $r--) {
    $moved = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r + 1, 1).Value = $moved
}
$ws.Cells.Item(14, 1).Value = "localdb"

# 4. Register the new named range for the category, matching the pattern
#    used by every other category name.
$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
